$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: human-readable labels
$ws.Range("A1").Value = "Edad"
$ws.Range("B1").Value = "Rama actividad, descripción"
$ws.Range("C1").Value = "Personas residentes viviendas familiares"
$ws.Range("D1").Value = "Rama actividad, código"
$ws.Range("E1").Value = "Aragón"
$ws.Range("F1").Value = "Sector actividad"
$ws.Range("G1").Value = "Sexo"

# Row 2: identifiers (measure/dimension urn-like codes)
$ws.Range("A2").Value = "iaest-measure:edad"
$ws.Range("B2").Value = "iaest-measure:rama-actividad-descripcion"
$ws.Range("C2").Value = "iaest-measure:personas-residentes-viviendas-familiares"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "iaest-measure:sector-actividad"
$ws.Range("G2").Value = "iaest-measure:sexo"

# Row 3: kind (medida / dim / null)
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "dim"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "medida"

# Row 4: data type
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:string"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("F4").Value = "xsd:string"
$ws.Range("G4").Value = "xsd:string"
